$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for columns B (TB), C (d2S), D (K), E (IP) for rows 2-15.
# Column F (Win) is unchanged; column G (sum) is recomputed as B+C+D+E.
$bcde = @(
        @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987),
        @(3.272327238179451, 1.626987699542094, 18.71679738969934, 13.86384647080068),
        @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 0.5333859586016987),
        @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987),
        @(0.04172184405617529, 0.3048912486333797, 3.223369029078222, 13.86384647080068),
        @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987),
        @(0.1169995834814548, 0.04103571897497393, 3.223369029078222, 0.5333859586016987),
        @(0.00009552326474482342, 0.002658071450198252, 0.1496068669990043, 0.5333859586016987),
        @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987),
        @(0.6545652718822623, 2919.202174992006, 0.7210945179870265, 13.86384647080068),
        @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987),
        @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 2797.565817734744),
        @(0.1169995834814548, 1.626987699542094, 0.1496068669990043, 0.5333859586016987),
        @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068)
)

$gsum = @(
        8.656069925401464,
        37.47995879822157,
        4.178645819794754,
        6.82939032824165,
        17.43382859256846,
        5.582307763322248,
        3.914790290136349,
        0.685746420315646,
        8.656069925401464,
        2934.441681252676,
        6.82939032824165,
        2803.186227190452,
        2.426980108624251,
        21.98653043760045
)

for ($i = 0; $i -lt 14; $i++) {
    $r = $i + 2
    $vals = $bcde[$i]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $gsum[$i]
}
